$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "59.832.92"
Set-TextValue "E2" "  +0.01%  "
Set-TextValue "D3" "2.306.53"
Set-TextValue "E3" "  -1.87%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "543.91"
Set-TextValue "E5" "  -0.09%  "
Set-TextValue "D6" "129.12"
Set-TextValue "E6" "  -2.28%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "0.574"
Set-TextValue "E8" "  -2.47%  "
Set-TextValue "D9" "2.303.65"
Set-TextValue "E9" "  -1.87%  "
Set-TextValue "E10" "  -0.10%  "
Set-TextValue "D11" "5.54"
Set-TextValue "E11" "  +0.64%  "
Set-TextValue "E12" "  -0.52%  "
Set-TextValue "D13" "0.333"
Set-TextValue "E13" "  +0.05%  "
Set-TextValue "D14" "23.30"
Set-TextValue "E14" "  -2.15%  "
Set-TextValue "B15" "WrappedliquidstakedEther2.0"
Set-TextValue "C15" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D15" "2.729.23"
Set-TextValue "E15" "  -1.45%  "
Set-TextValue "B16" "WrappedBTC"
Set-TextValue "C16" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D16" "59.929.81"
Set-TextValue "E16" "  +0.08%  "
Set-TextValue "E17" "  +0.18%  "
Set-TextValue "D18" "2.320.09"
Set-TextValue "E18" "  -1.31%  "
Set-TextValue "D19" "10.51"
Set-TextValue "E19" "  -1.38%  "
Set-TextValue "B20" "BitcoinCash"
Set-TextValue "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "311.90"
Set-TextValue "E20" "  -0.38%  "
Set-TextValue "B21" "Polkadot"
Set-TextValue "C21" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D21" "4.05"
Set-TextValue "E21" "  -2.54%  "
Set-TextValue "D22" "6.54"
Set-TextValue "E22" "  -4.06%  "
Set-TextValue "E23" "  +0.19%  "
Set-TextValue "D24" "63.86"
Set-TextValue "E24" "  +1.29%  "
Set-TextValue "D25" "0.169"
Set-TextValue "E25" "  -0.82%  "
Set-TextValue "D27" "7.78"
Set-TextValue "E27" "  -1.34%  "
Set-TextValue "D28" "1.36"
Set-TextValue "E28" "  +2.37%  "
Set-TextValue "E29" "  +7.69%  "
Set-TextValue "D30" "171.29"
Set-TextValue "E30" "  -0.05%  "
Set-TextValue "E31" "  -0.89%  "
Set-TextValue "D32" "0.0₃0721"
Set-TextValue "E32" "  -0.58%  "
Set-TextValue "D33" "5.92"
Set-TextValue "E33" "  +0.21%  "
Set-TextValue "D34" "0.381"
Set-TextValue "E34" "  +0.25%  "
Set-TextValue "E35" "  -4.50%  "
Set-TextValue "B36" "USDe"
Set-TextValue "C36" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  +0.02%  "
Set-TextValue "B37" "EthereumClassic"
Set-TextValue "C37" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D37" "17.88"
Set-TextValue "E37" "  -0.60%  "
Set-TextValue "E38" "  +0.01%  "
Set-TextValue "D39" "4.06"
Set-TextValue "E39" "  -1.65%  "
Set-TextValue "D40" "314.80"
Set-TextValue "E40" "  +0.59%  "
Set-TextValue "D41" "37.89"
Set-TextValue "E41" "  -0.40%  "
Set-TextValue "D42" "1.51"
Set-TextValue "E42" "  -0.92%  "
Set-TextValue "D43" "136.51"
Set-TextValue "E43" "  -3.98%  "
Set-TextValue "D44" "3.47"
Set-TextValue "E44" "  +0.59%  "
Set-TextValue "D45" "0.0936"
Set-TextValue "E45" "  -1.68%  "
Set-TextValue "D46" "18.95"
Set-TextValue "E46" "  -1.57%  "
Set-TextValue "D47" "0.561"
Set-TextValue "E47" "  -0.07%  "
Set-TextValue "D48" "0.0493"
Set-TextValue "E48" "  -0.57%  "
Set-TextValue "D49" "0.0214"
Set-TextValue "E49" "  +0.18%  "
Set-TextValue "D50" "0.0₆0213"
Set-TextValue "E50" "  +3.50%  "
Set-TextValue "D51" "16.77"
Set-TextValue "E51" "  -0.78%  "
